# "fixed grade for styling"
#
# Sheet1 ("SRS Grade Rubric"):
#  - B14:C14 had stray 0 values entered for a purely-styling/separator row;
#    clear them back out so the cells are blank again (style stays s="6").
#  - B15 ("Total" for that section) bumps from 7 to 8.
#  - B71 ("Total Grade") bumps from 97.5 to 98.5.
#  - Leave the view scrolled down near the bottom of the rubric, with B71
#    (the corrected total) selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the leftover 0 values in B14:C14, keeping their existing style.
$ws.Range("B14:C14").ClearContents()

# Corrected totals.
$ws.Range("B15").Value = 8
$ws.Range("B71").Value = 98.5

# Scroll the window down and leave B71 selected, matching the saved view.
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B71").Select()
